$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = "27.820.10"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "1.768.83"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.83%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.00"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("E7").Value = "  -7.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.42"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.129"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07470"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.008"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.01"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.194"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.328"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "1.769.70"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001072"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06557"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.27"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.009"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.10"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.166"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").Value = "27.886.49"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.36"
$ws.Range("E24").Value = "  -4.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.406"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.28"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.387"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.77"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "1.974.26"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.312"
$ws.Range("E30").Value = "  +2.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.82"
$ws.Range("E31").Value = "  -3.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.977"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.726"
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09055"
$ws.Range("E34").Value = "  -4.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.38"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2203"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06220"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6579"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.191"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.430"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.026"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.88"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.798"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5982"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.976"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06930"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.121"
$ws.Range("E51").Value = "  -4.11%  "

# --- Rows 39/40: VeChain and InternetComputer(DFINITY) swapped order ---
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02291"
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.123"
$ws.Range("E40").Value = "  -1.47%  "
